$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.189.93'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '2.231.60'
$ws.Range("E3").Value = '  -0.13%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = "'259.10"
$ws.Range("E5").Value = '  +2.92%  '
$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = '  +1.73%  '
$ws.Range("D7").Value = "'78.32"
$ws.Range("E7").Value = '  +4.13%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = '  +0.90%  '
$ws.Range("D10").Value = "'43.19"
$ws.Range("E10").Value = '  +5.25%  '
$ws.Range("D11").Value = "'0.0921"
$ws.Range("E11").Value = '  -0.13%  '
$ws.Range("D12").Value = "'7.05"
$ws.Range("E12").Value = '  +2.72%  '
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("D14").Value = '2.561.74'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("D15").Value = "'14.61"
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '2.242.99'
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = "'0.792"
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").Value = '43.146.56'
$ws.Range("E18").Value = '  +0.54%  '
$ws.Range("E19").Value = '  +0.42%  '
$ws.Range("D20").Value = "'71.32"
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D21").Value = "'6.02"
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("D22").Value = "'2.33"
$ws.Range("E22").Value = '  +6.15%  '
$ws.Range("D23").Value = "'231.15"
$ws.Range("E23").Value = '  +0.65%  '
$ws.Range("D24").Value = "'9.27"
$ws.Range("E24").Value = '  -1.42%  '
$ws.Range("E25").Value = '  -0.10%  '
$ws.Range("D26").Value = "'42.73"
$ws.Range("E26").Value = '  +9.08%  '
$ws.Range("D27").Value = "'10.85"
$ws.Range("E27").Value = '  +1.47%  '
$ws.Range("E28").Value = '  -2.28%  '
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = '  +0.47%  '
$ws.Range("D30").Value = "'2.20"
$ws.Range("E30").Value = '  +2.55%  '
$ws.Range("D31").Value = "'173.52"
$ws.Range("E31").Value = '  +1.54%  '
$ws.Range("D32").Value = "'20.46"
$ws.Range("E32").Value = '  +1.42%  '
$ws.Range("D33").Value = "'0.0873"
$ws.Range("E33").Value = '  +10.07%  '
$ws.Range("D34").Value = "'5.27"
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("D35").Value = "'0.123"
$ws.Range("E35").Value = '  +1.22%  '
$ws.Range("D36").Value = "'0.0368"
$ws.Range("E36").Value = '  +12.50%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = "'0.108"
$ws.Range("E37").Value = '  -3.42%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = "'4.46"
$ws.Range("E38").Value = '  -0.59%  '
$ws.Range("D39").Value = "'13.34"
$ws.Range("E39").Value = '  +8.16%  '
$ws.Range("D40").Value = "'2.91"
$ws.Range("E40").Value = '  +19.86%  '
$ws.Range("D41").Value = "'2.14"
$ws.Range("E41").Value = '  +1.97%  '
$ws.Range("B42").Value = 'MultiversX'
$ws.Range("C42").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D42").Value = "'61.88"
$ws.Range("E42").Value = '  +5.86%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = "'0.204"
$ws.Range("E43").Value = '  -0.57%  '
$ws.Range("D44").Value = "'5.33"
$ws.Range("E44").Value = '  -0.66%  '
$ws.Range("D45").Value = "'104.04"
$ws.Range("E45").Value = '  +0.93%  '
$ws.Range("D46").Value = "'8.55"
$ws.Range("E46").Value = '  -1.04%  '
$ws.Range("D47").Value = "'0.477"
$ws.Range("E47").Value = '  -2.24%  '
$ws.Range("D48").Value = "'0.0980"
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("E49").Value = '  +0.70%  '
$ws.Range("E50").Value = '  +0.84%  '
$ws.Range("E51").Value = '  +22.76%  '
